$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "full data for NY sites" — fill in the remaining fields (event image URL
# and venue lat/long) for both rows, extending the used range from N to Q.

# Row 1: The Kids in the Hall @ The Town Hall
$ws.Range("O1").Value = "http://upload.wikimedia.org/wikipedia/commons/8/8e/ETalk2008-Kids_In_The_Hall.jpg"
$ws.Range("P1").Value = 40.755986
$ws.Range("Q1").Value = -73.984712000000002

# Row 2: Carol Burnett @ St. George Theatre
$ws.Range("O2").Value = "http://upload.wikimedia.org/wikipedia/commons/8/83/Carol_Burnett_-_1974.jpg"
$ws.Range("P2").Value = 40.643332999999998
$ws.Range("Q2").Value = -74.078889000000004

# Scroll the view toward the newly-populated columns and leave the new
# bottom-right cell selected, matching the author's final cursor position.
try {
    $excel.ActiveWindow.ScrollColumn = 10
} catch {
}

$ws.Range("P2").Select()
